# Bug fix, everything changed to long
# Update benchmark result cells (Zeit / Speicher / Zugriffe / Vergleiche)
# across all 9 algorithm-result worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("InversTeilsortiert1000.dat")
$ws.Range("B2").Value = 5339600
$ws.Range("B3").Value = 5461700
$ws.Range("C3").Value = 32104
$ws.Range("B4").Value = 1178700
$ws.Range("B5").Value = 57000
$ws.Range("C5").Value = 32128
$ws.Range("B6").Value = 2150100
$ws.Range("C6").Value = 32128
$ws.Range("B7").Value = 1406100
$ws.Range("D7").Value = 19569
$ws.Range("E7").Value = 4523
$ws.Range("B8").Value = 4206600

$ws = $wb.Worksheets.Item("InversTeilsortiert10000.dat")
$ws.Range("B2").Value = 101770100
$ws.Range("B3").Value = 79061100
$ws.Range("C3").Value = 320104
$ws.Range("B4").Value = 2994500
$ws.Range("B5").Value = 521900
$ws.Range("C5").Value = 320128
$ws.Range("B6").Value = 83063300
$ws.Range("C6").Value = 320128
$ws.Range("B7").Value = 2184800
$ws.Range("D7").Value = 339711
$ws.Range("E7").Value = 91237
$ws.Range("B8").Value = 60010300

$ws = $wb.Worksheets.Item("InversTeilsortiert100000.dat")
$ws.Range("B2").Value = 9363035700
$ws.Range("B3").Value = 8813514900
$ws.Range("C3").Value = 3200104
$ws.Range("B4").Value = 10538400
$ws.Range("B5").Value = 3946400
$ws.Range("C5").Value = 3200128
$ws.Range("B6").Value = 2253614200
$ws.Range("C6").Value = 3200128
$ws.Range("B7").Value = 14715000
$ws.Range("D7").Value = 3753264
$ws.Range("E7").Value = 1029088
$ws.Range("B8").Value = 4301001600

$ws = $wb.Worksheets.Item("Random1000.dat")
$ws.Range("B2").Value = 622300
$ws.Range("B3").Value = 860800
$ws.Range("C3").Value = 32104
$ws.Range("B4").Value = 79300
$ws.Range("B5").Value = 3700
$ws.Range("C5").Value = 32128
$ws.Range("B6").Value = 143700
$ws.Range("C6").Value = 32128
$ws.Range("B7").Value = 101200
$ws.Range("D7").Value = 3773814
$ws.Range("E7").Value = 1033938
$ws.Range("B8").Value = 543200

$ws = $wb.Worksheets.Item("Random10000.dat")
$ws.Range("B2").Value = 47413000
$ws.Range("B3").Value = 77494700
$ws.Range("C3").Value = 320104
$ws.Range("B4").Value = 1112200
$ws.Range("B5").Value = 102200
$ws.Range("C5").Value = 320128
$ws.Range("B6").Value = 15084900
$ws.Range("C6").Value = 320128
$ws.Range("B7").Value = 1204000
$ws.Range("D7").Value = 4102989
$ws.Range("E7").Value = 1123663
$ws.Range("B8").Value = 51445200

$ws = $wb.Worksheets.Item("Random100000.dat")
$ws.Range("B2").Value = 14745218500
$ws.Range("B3").Value = 11630467800
$ws.Range("C3").Value = 3200104
$ws.Range("B4").Value = 12822200
$ws.Range("B5").Value = 90700
$ws.Range("C5").Value = 3200128
$ws.Range("B6").Value = 1664846300
$ws.Range("C6").Value = 3200128
$ws.Range("B7").Value = 14056600
$ws.Range("D7").Value = 7639881
$ws.Range("E7").Value = 2102627
$ws.Range("B8").Value = 4713854200

$ws = $wb.Worksheets.Item("Teilsortiert1000.dat")
$ws.Range("B2").Value = 433500
$ws.Range("B3").Value = 1092400
$ws.Range("C3").Value = 32104
$ws.Range("B4").Value = 88000
$ws.Range("B5").Value = 1400
$ws.Range("C5").Value = 32128
$ws.Range("B6").Value = 4343600
$ws.Range("C6").Value = 32128
$ws.Range("B7").Value = 116100
$ws.Range("D7").Value = 7663869
$ws.Range("E7").Value = 2108623
$ws.Range("B8").Value = 556300

$ws = $wb.Worksheets.Item("Teilsortiert10000.dat")
$ws.Range("B2").Value = 36767200
$ws.Range("B3").Value = 18760900
$ws.Range("C3").Value = 320104
$ws.Range("B4").Value = 1065500
$ws.Range("B5").Value = 13100
$ws.Range("C5").Value = 320128
$ws.Range("B6").Value = 21356100
$ws.Range("C6").Value = 320128
$ws.Range("B7").Value = 965900
$ws.Range("D7").Value = 7958565
$ws.Range("E7").Value = 2186855
$ws.Range("B8").Value = 45989000

$ws = $wb.Worksheets.Item("Teilsortiert100000.dat")
$ws.Range("B2").Value = 6309057400
$ws.Range("B3").Value = 2495298900
$ws.Range("C3").Value = 3200104
$ws.Range("B4").Value = 8964200
$ws.Range("B5").Value = 94200
$ws.Range("C5").Value = 3200128
$ws.Range("B6").Value = 2198457100
$ws.Range("C6").Value = 3200128
$ws.Range("B7").Value = 10466200
$ws.Range("D7").Value = 12150057
$ws.Range("E7").Value = 3384019
$ws.Range("B8").Value = 4325120600
